$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 12
$ws.Range("D12").Value = "“(번역개정판) 파이썬 라이브러리를 활용한 머신러닝” 주피터 노트북 업데이트 및 에러타 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/03/08/%eb%b2%88%ec%97%ad%ea%b0%9c%ec%a0%95%ed%8c%90-%ed%8c%8c%ec%9d%b4%ec%8d%ac-%eb%9d%bc%ec%9d%b4%eb%b8%8c%eb%9f%ac%eb%a6%ac%eb%a5%bc-%ed%99%9c%ec%9a%a9%ed%95%9c-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d/"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Interpretable Convolutional Neural Networks"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1447&mod=document&pageid=1"

# Row 39
$ws.Range("D39").Value = "Probability concepts explained: Bayesian inference for parameter estimation"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-Bayesian-inference-for-parameter-estimation-1"

# Row 41
$ws.Range("D41").Value = "Service Mesh 에서의 Control Plane, Consul"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/service-mesh-%ec%97%90%ec%84%9c%ec%9d%98-control-plane-consul/"

# Row 51
$ws.Range("D51").Value = "[세이버메트릭스] 타자의 타율, OPS, 삼진비율, 볼넷비율 중 어떤 것을 가장 신뢰할 수 있을까?"
$ws.Range("E51").Value = "https://bskyvision.com/1134"
